$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-02-12T16:16:10.958439"

# Row 2 (Current_Ct_Day, Current_Pct_Ct, Current_Ct_Tokens, Current_Pct_Tokens, Last_Reset)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = [double]"6.944444444444444e-05"
$ws.Range("K2").Value = 472
$ws.Range("L2").Value = 0.000944
$ws.Range("M2").Value = $timestamp

# Rows 3-14: counters reset to 0, Last_Reset refreshed
for ($row = 3; $row -le 14; $row++) {
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = $timestamp
}
